# Applies the cryptos-list refresh described in the commit diff.
# Cells in columns D/E store text (e.g. "1.00", "  -0.52%  "), not numbers,
# so numeric-looking values are written with a temporary Text number format
# to stop Excel from auto-coercing them into the Number type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.906.17"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "2.296.82"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.10"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.12"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -3.69%  "
$ws.Range("E7").Value = "  -0.50%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.497"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -3.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.46"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -4.78%  "
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("E12").Value = "  -4.08%  "
$ws.Range("E13").Value = "  +1.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.78"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +6.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.78"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("D16").Value = "2.651.30"
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("D17").Value = "2.288.77"
$ws.Range("E17").Value = "  -1.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.808"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("D19").Value = "42.783.93"
$ws.Range("E19").Value = "  -0.80%  "
$ws.Range("D20").Value = "0.0₃0900"
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.54"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -2.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.03"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -1.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.51"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.29"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -0.83%  "
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  -3.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.32"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -3.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.98"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.77"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -2.88%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.04"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -6.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.13"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -1.19%  "
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.68"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +3.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.93"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -2.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "16.97"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0692"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -2.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.102"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.82"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -2.62%  "
$ws.Range("E41").Value = "  -3.90%  "
$ws.Range("E42").Value = "  -1.21%  "
$ws.Range("E43").Value = "  -3.41%  "
$ws.Range("D44").Value = "1.991.44"
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("E45").Value = "  -1.65%  "
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.52"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -6.29%  "
$ws.Range("E48").Value = "  -3.62%  "
$ws.Range("D49").Value = "2.520.79"
$ws.Range("E49").Value = "  -0.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.12"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -3.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.57"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -5.61%  "
